$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that immediately follows
#    the title (Heading1) paragraph.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play Dragon's Fire Free Slot Game |
#    Stunning 2D and 3D Animations") right before the paragraph that
#    currently holds the "Prompt: ..." text (the last paragraph in the
#    body). Inserting at (End - 1) of the paragraph immediately
#    preceding it reliably produces a clean, separate new paragraph
#    (with the exact <w:r/><w:r><w:rPr>...</w:rPr>...</w:r> run shape)
#    without disturbing either neighbour.
# ---------------------------------------------------------------------
$precedingPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$insertPos = $precedingPara.Range.End - 1
$insertRng = $d.Range($insertPos, $insertPos)
$insertRng.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dragon&#39;s Fire Free Slot Game | Stunning 2D and 3D Animations</w:t></w:r></w:p>")

# ---------------------------------------------------------------------
# 3. Replace the text of the former "Prompt: ..." paragraph (now the
#    last paragraph again) with the new meta-description sentence,
#    keeping its existing italic run formatting intact.
# ---------------------------------------------------------------------
$promptPara = $d.Paragraphs($d.Paragraphs.Count)
$pr = $promptPara.Range
$textRng = $d.Range($pr.Start, $pr.End - 1)
$textRng.Text = "Experience Dragon's Fire, a stunning online slot game with 2D and 3D animations. Play for free on desktops and mobile devices."

Write-Output "done"
